# Scheduled runner update: refresh computed profit/cost figures across
# all item sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) in the Ixion
# Profits workbook. Values below mirror the latest market-data snapshot.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 26429.285
$ws.Range("J10").Value = 26429.285
$ws.Range("L10").Value = 26429.285
$ws.Range("N10").Value = -27015.285
$ws.Range("H40").Value = 2010.3066
$ws.Range("I40").Value = 2006.8096
$ws.Range("J40").Value = 2028.6666
$ws.Range("K40").Value = 2006.8096
$ws.Range("L40").Value = 2028.6666
$ws.Range("M40").Value = -1831.8096
$ws.Range("N40").Value = -2378.6666
$ws.Range("H95").Value = 43000
$ws.Range("J95").Value = 43000
$ws.Range("L95").Value = 43000
$ws.Range("N95").Value = -48492
$ws.Range("H97").Value = 928.5714
$ws.Range("I97").Value = 600
$ws.Range("J97").Value = 1060
$ws.Range("K97").Value = 1800
$ws.Range("L97").Value = 3180
$ws.Range("M97").Value = -1304
$ws.Range("N97").Value = -4172
$ws.Range("H132").Value = 2855.5278
$ws.Range("I132").Value = 2077.16
$ws.Range("J132").Value = 4624.5454
$ws.Range("K132").Value = 6231.48
$ws.Range("L132").Value = 13873.6362
$ws.Range("M132").Value = -3701.48
$ws.Range("N132").Value = -18933.6362

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 4000
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H61").Value = 318508.34
$ws.Range("I61").Value = 7402.3335
$ws.Range("J61").Value = 912438
$ws.Range("K61").Value = 7402.3335
$ws.Range("L61").Value = 912438
$ws.Range("M61").Value = -7190.3335
$ws.Range("N61").Value = -912862
$ws.Range("H98").Value = 43000
$ws.Range("J98").Value = 43000
$ws.Range("L98").Value = 43000
$ws.Range("N98").Value = -48990
$ws.Range("H101").Value = 59801
$ws.Range("J101").Value = 59801
$ws.Range("L101").Value = 59801
$ws.Range("N101").Value = -66291
$ws.Range("H104").Value = 40862.5
$ws.Range("J104").Value = 40862.5
$ws.Range("L104").Value = 40862.5
$ws.Range("N104").Value = -47850.5
$ws.Range("H132").Value = 2944795.5
$ws.Range("I132").Value = 2733.4
$ws.Range("J132").Value = 7147741.5
$ws.Range("K132").Value = 8200.200000000001
$ws.Range("L132").Value = 21443224.5
$ws.Range("M132").Value = -5670.200000000001
$ws.Range("N132").Value = -21448284.5
$ws.Range("H136").Value = 318508.34
$ws.Range("I136").Value = 7402.3335
$ws.Range("J136").Value = 912438
$ws.Range("K136").Value = 22207.0005
$ws.Range("L136").Value = 2737314
$ws.Range("M136").Value = -19657.0005
$ws.Range("N136").Value = -2742414

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 13157
$ws.Range("J76").Value = 13157
$ws.Range("L76").Value = 13157
$ws.Range("N76").Value = -13787
$ws.Range("H79").Value = 13157
$ws.Range("J79").Value = 13157
$ws.Range("L79").Value = 13157
$ws.Range("N79").Value = -15341
$ws.Range("H134").Value = 44960.117
$ws.Range("I134").Value = 9774
$ws.Range("J134").Value = 86010.586
$ws.Range("K134").Value = 29322
$ws.Range("L134").Value = 258031.758
$ws.Range("M134").Value = -26787
$ws.Range("N134").Value = -263101.758

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 1250
$ws.Range("I14").Value = 500
$ws.Range("J14").Value = 2000
$ws.Range("K14").Value = 500
$ws.Range("L14").Value = 2000
$ws.Range("M14").Value = -330
$ws.Range("N14").Value = -2340
$ws.Range("H31").Value = 4694.4746
$ws.Range("I31").Value = 2124.3242
$ws.Range("J31").Value = 9017
$ws.Range("K31").Value = 2124.3242
$ws.Range("L31").Value = 9017
$ws.Range("M31").Value = -1829.3242
$ws.Range("N31").Value = -9607
$ws.Range("H34").Value = 4694.4746
$ws.Range("I34").Value = 2124.3242
$ws.Range("J34").Value = 9017
$ws.Range("K34").Value = 2124.3242
$ws.Range("L34").Value = 9017
$ws.Range("M34").Value = -1922.3242
$ws.Range("N34").Value = -9421
$ws.Range("H99").Value = 11208.333
$ws.Range("I99").Value = 15900
$ws.Range("J99").Value = 1825
$ws.Range("K99").Value = 15900
$ws.Range("L99").Value = 1825
$ws.Range("M99").Value = -14402
$ws.Range("N99").Value = -4821
$ws.Range("H122").Value = 5456.4
$ws.Range("I122").Value = 4320.5
$ws.Range("K122").Value = 12961.5
$ws.Range("M122").Value = -10511.5
$ws.Range("H126").Value = 11208.333
$ws.Range("I126").Value = 15900
$ws.Range("J126").Value = 1825
$ws.Range("K126").Value = 47700
$ws.Range("L126").Value = 5475
$ws.Range("M126").Value = -45230
$ws.Range("N126").Value = -10415
$ws.Range("H132").Value = 2445.3
$ws.Range("I132").Value = 1247.4
$ws.Range("K132").Value = 3742.2
$ws.Range("M132").Value = -1212.2
$ws.Range("H134").Value = 195970.25
$ws.Range("I134").Value = 3838.457
$ws.Range("J134").Value = 591535.7
$ws.Range("K134").Value = 11515.371
$ws.Range("L134").Value = 1774607.1
$ws.Range("M134").Value = -8980.370999999999
$ws.Range("N134").Value = -1779677.1

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1538971.1
$ws.Range("I113").Value = 2000499.2
$ws.Range("J113").Value = 714813.7
$ws.Range("K113").Value = 6001497.6
$ws.Range("L113").Value = 2144441.1
$ws.Range("M113").Value = -5999327.6
$ws.Range("N113").Value = -2148781.1
$ws.Range("H131").Value = 1667688.1
$ws.Range("I131").Value = 11111856
$ws.Range("J131").Value = 1070.3334
$ws.Range("K131").Value = 33335568
$ws.Range("L131").Value = 3211.0002
$ws.Range("M131").Value = -33330528
$ws.Range("N131").Value = -13291.0002
$ws.Range("H132").Value = 2129.7856
$ws.Range("I132").Value = 548.8333
$ws.Range("J132").Value = 3315.5
$ws.Range("K132").Value = 4939.4997
$ws.Range("L132").Value = 29839.5
$ws.Range("M132").Value = -2409.4997
$ws.Range("N132").Value = -34899.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2290
$ws.Range("I102").Value = 2141.5386
$ws.Range("J102").Value = 2933.3333
$ws.Range("K102").Value = 2141.5386
$ws.Range("L102").Value = 2933.3333
$ws.Range("M102").Value = -519.5385999999999
$ws.Range("N102").Value = -6177.3333
$ws.Range("H104").Value = 20671
$ws.Range("J104").Value = 20671
$ws.Range("L104").Value = 20671
$ws.Range("N104").Value = -27659
$ws.Range("H121").Value = 30000
$ws.Range("J121").Value = 30000
$ws.Range("L121").Value = 30000
$ws.Range("N121").Value = -33494
$ws.Range("H126").Value = 25812.5
$ws.Range("I126").Value = 33250
$ws.Range("K126").Value = 99750
$ws.Range("M126").Value = -97280
$ws.Range("H132").Value = 5061.6733
$ws.Range("I132").Value = 9165.277
$ws.Range("J132").Value = 2678.9355
$ws.Range("K132").Value = 27495.831
$ws.Range("L132").Value = 8036.806500000001
$ws.Range("M132").Value = -24965.831
$ws.Range("N132").Value = -13096.8065
$ws.Range("H141").Value = 64356.25
$ws.Range("J141").Value = 64356.25
$ws.Range("L141").Value = 64356.25
$ws.Range("N141").Value = -74716.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 6375
$ws.Range("J11").Value = 6375
$ws.Range("L11").Value = 6375
$ws.Range("N11").Value = -6655
$ws.Range("H40").Value = 76925970
$ws.Range("I40").Value = 90912056
$ws.Range("K40").Value = 90912056
$ws.Range("M40").Value = -90911920
$ws.Range("H46").Value = 976.9091
$ws.Range("I46").Value = 550.3333
$ws.Range("J46").Value = 1488.8
$ws.Range("K46").Value = 550.3333
$ws.Range("L46").Value = 1488.8
$ws.Range("M46").Value = -362.3333
$ws.Range("N46").Value = -1864.8
$ws.Range("H64").Value = 41600
$ws.Range("J64").Value = 41600
$ws.Range("L64").Value = 41600
$ws.Range("N64").Value = -42050
$ws.Range("H67").Value = 41600
$ws.Range("J67").Value = 41600
$ws.Range("L67").Value = 41600
$ws.Range("N67").Value = -43160
$ws.Range("H101").Value = 31181
$ws.Range("J101").Value = 31181
$ws.Range("L101").Value = 31181
$ws.Range("N101").Value = -37671
$ws.Range("H106").Value = 34500
$ws.Range("J106").Value = 34500
$ws.Range("L106").Value = 34500
$ws.Range("N106").Value = -37024
$ws.Range("H122").Value = 7946352
$ws.Range("I122").Value = 7946352
$ws.Range("K122").Value = 23839056
$ws.Range("M122").Value = -23836606

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 20307.8
$ws.Range("J92").Value = 20307.8
$ws.Range("L92").Value = 20307.8
$ws.Range("N92").Value = -25299.8
$ws.Range("H97").Value = 42628.6
$ws.Range("J97").Value = 42628.6
$ws.Range("L97").Value = 42628.6
$ws.Range("N97").Value = -44610.6
$ws.Range("H104").Value = 43000
$ws.Range("J104").Value = 43000
$ws.Range("L104").Value = 43000
$ws.Range("N104").Value = -49988
$ws.Range("H107").Value = 43478850
$ws.Range("I107").Value = 71429040
$ws.Range("J107").Value = 774.7778
$ws.Range("K107").Value = 214287120
$ws.Range("L107").Value = 2324.3334
$ws.Range("M107").Value = -214285200
$ws.Range("N107").Value = -6164.3334
$ws.Range("H122").Value = 1918.75
$ws.Range("I122").Value = 1861.5
$ws.Range("J122").Value = 2205
$ws.Range("K122").Value = 5584.5
$ws.Range("L122").Value = 6615
$ws.Range("M122").Value = -3134.5
$ws.Range("N122").Value = -11515
$ws.Range("H140").Value = 45721.617
$ws.Range("J140").Value = 45721.617
$ws.Range("L140").Value = 45721.617
$ws.Range("N140").Value = -56081.617
